$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new "Person Responsible" column before the existing "Date Last
# Edited" column F, shifting old F ("Date Last Edited") -> G and old G
# ("Comments") -> H. Copy+Insert carries over the formatting (fill/font/
# alignment) so every shifted column keeps its exact original styling, and
# the new column starts out looking like column F too.
$fWidth = $ws.Columns.Item(6).ColumnWidth
$ws.Columns.Item(6).Copy()
$ws.Columns.Item(6).Insert()
$ws.Columns.Item(6).ColumnWidth = $fWidth

# Bump the format version
$ws.Range("C2").Value = "v0.2.0"

# Fill in the new "Person Responsible" column (now column F)
$ws.Range("F3").Value = "Person Responsible"
$ws.Range("F4").Value = "personResponsible"
$ws.Range("F5").Value = "Person responsible that added this source and the corresponding entries"
$ws.Range("F6").Value = "-"
$ws.Range("F7").Value = "[text]"
$ws.Range("F8").Value = "Tester"
$ws.Range("F9").Value = "Tester"
